{"js": "const body = context.document.body;\n\n// --- Edit 1 -----------------------------------------------------------\n// \"Familywise error rates and hidden multiplicity in ANOVA \" ->\n// \"Understanding \" + italic \"p\" + \"-hacking via multiple testing \"\n// (the trailing \"[assignment? - p hacking]\" run is left untouched)\nconst results1 = body.search(\n  \"Familywise error rates and hidden multiplicity in ANOVA \",\n  { matchCase: true }\n);\nresults1.load(\"items\");\nawait context.sync();\n\nif (results1.items.length > 0) {\n  const r1 = results1.items[0];\n\n  // Insert the replacement pieces immediately after the matched range,\n  // each insertion lands right after r1 (pushing earlier insertions\n  // further along), so insert in reverse order.\n  const seg4 = r1.insertText(\n    \"-hacking via multiple testing \",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n\n  const seg3 = r1.insertText(\"p\", Word.InsertLocation.after);\n  seg3.font.italic = true;\n  await context.sync();\n\n  const seg2 = r1.insertText(\" \", Word.InsertLocation.after);\n  seg2.font.italic = true;\n  await context.sync();\n\n  const seg1 = r1.insertText(\"Understanding\", Word.InsertLocation.after);\n  await context.sync();\n\n  // Remove the original matched text now that the replacement is in place.\n  r1.delete();\n  await context.sync();\n}\n\n// --- Edit 2 -------------------------------------------------------------\n// \"Understanding Confidence Intervals via sequential testing \" ->\n// \"Understanding the relationships between Confidence Intervals, \" +\n// italic \"p\" + \"-values, and statistical power\"\nconst results2 = body.search(\n  \"Understanding Confidence Intervals via sequential testing \",\n  { matchCase: true }\n);\nresults2.load(\"items\");\nawait context.sync();\n\nif (results2.items.length > 0) {\n  const r2 = results2.items[0];\n\n  const t6 = r2.insertText(\n    \"-values, and statistical power\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n\n  const t5 = r2.insertText(\"p\", Word.InsertLocation.after);\n  t5.font.italic = true;\n  await context.sync();\n\n  const t4 = r2.insertText(\", \", Word.InsertLocation.after);\n  await context.sync();\n\n  const t3 = r2.insertText(\"Confidence Intervals\", Word.InsertLocation.after);\n  await context.sync();\n\n  const t2 = r2.insertText(\n    \"the relationships between \",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n\n  const t1 = r2.insertText(\"Understanding \", Word.InsertLocation.after);\n  await context.sync();\n\n  r2.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction InsertAfterRange($baseRange, $text) {\n    $r = $baseRange.Duplicate\n    $r.Collapse(0)  # wdCollapseEnd\n    $r.InsertAfter($text)\n    return $r\n}\n\n# --- Edit 1 --------------------------------------------------------------\n# \"Familywise error rates and hidden multiplicity in ANOVA \" ->\n# \"Understanding \" + italic \"p\" + \"-hacking via multiple testing \"\n# (the trailing \"[assignment? - p hacking]\" text is left untouched)\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.Text = \"Familywise error rates and hidden multiplicity in ANOVA \"\n$found1 = $find1.Execute()\n\nif ($found1) {\n    $seg4 = InsertAfterRange $rng1 \"-hacking via multiple testing \"\n\n    $seg3 = InsertAfterRange $rng1 \"p\"\n    $seg3.Italic = 1\n\n    $seg2 = InsertAfterRange $rng1 \" \"\n    $seg2.Italic = 1\n\n    $seg1 = InsertAfterRange $rng1 \"Understanding\"\n\n    $rng1.Delete()\n}\n\n# --- Edit 2 ----------------------------------------------------------------\n# \"Understanding Confidence Intervals via sequential testing \" ->\n# \"Understanding the relationships between Confidence Intervals, \" +\n# italic \"p\" + \"-values, and statistical power\"\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.Text = \"Understanding Confidence Intervals via sequential testing \"\n$found2 = $find2.Execute()\n\nif ($found2) {\n    $t6 = InsertAfterRange $rng2 \"-values, and statistical power\"\n\n    $t5 = InsertAfterRange $rng2 \"p\"\n    $t5.Italic = 1\n\n    $t4 = InsertAfterRange $rng2 \", \"\n\n    $t3 = InsertAfterRange $rng2 \"Confidence Intervals\"\n\n    $t2 = InsertAfterRange $rng2 \"the relationships between \"\n\n    $t1 = InsertAfterRange $rng2 \"Understanding \"\n\n    $rng2.Delete()\n}\n\nWrite-Output \"done\"\n"}
